$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): rename the id-builder columns and append "-JKO" to every
# course name; also insert a new "Violence Response (1 hr)" course before the
# trailing "RandomCourse" column.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "dodid"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "cleanName"
$ws.Range("E1").Value = "fullName"
$ws.Range("F1").Value = "courseName"
$ws.Range("G1").Value = "DHA Accommodations (1 hr)-JKO"
$ws.Range("H1").Value = "Leadership Training (4 hrs)-JKO"
$ws.Range("I1").Value = "MHS Customer Service (1 hr)-JKO"
$ws.Range("J1").Value = "Counterintelligence (1 hr)-JKO"
$ws.Range("K1").Value = "HIPAA Training (1 hr)-JKO"
$ws.Range("L1").Value = "Supervisor Safety Training (2 hrs)-JKO"
# M1:O1 are brand-new header cells - pick up the same bold/border/centred
# style already used by the rest of row 1 (xf index 1) by copying formats
# from the preceding header cell instead of re-deriving a new style.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Employee Safety (1 hr)-JKO"
$ws.Range("N1").Value = "Violence Response (1 hr)-JKO"
$ws.Range("O1").Value = "RandomCourse-JKO"

# ---------------------------------------------------------------------------
# Row 2 - John Doe
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "johndoe"
$ws.Range("E2").Value = "John Doe"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = 45536
$ws.Range("G2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H2").Value = 45536
$ws.Range("H2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I2").Value = 45536
$ws.Range("I2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""

# ---------------------------------------------------------------------------
# Row 3 - Andrew Hartmann
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "andrewhartmann"
$ws.Range("E3").Value = "Andrew Hartmann"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = 45536
$ws.Range("J3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K3").Value = 45536
$ws.Range("K3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L3").Value = 45505
$ws.Range("L3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M3").Value = 45505
$ws.Range("M3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""

# ---------------------------------------------------------------------------
# Row 4 - Nick Fletcher
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "nickfletcher"
$ws.Range("E4").Value = "Nick Fletcher"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = 45536
$ws.Range("J4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K4").Value = 45536
$ws.Range("K4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L4").Value = 45505
$ws.Range("L4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("M4").Value = 45505
$ws.Range("M4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = 36527
$ws.Range("O4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# Row 5 - John Cena
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "johncena"
$ws.Range("E5").Value = "John Cena"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = 45536
$ws.Range("G5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H5").Value = 45536
$ws.Range("H5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I5").Value = 45536
$ws.Range("I5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = 45536
$ws.Range("N5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("O5").Value = 36161
$ws.Range("O5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
